$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove rows 9 and 10 (data for these commesse no longer present)
$ws.Range("A10").EntireRow.Delete()
$ws.Range("A9").EntireRow.Delete()

$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B3").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("L3").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L3").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("A6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A2").Value = 233333
$ws.Range("B2").Value = 45903
$ws.Range("C2").Value = "STAMPATO"
$ws.Range("D2").Value = "CAMPO VUOTO"
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 12
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = "1"
$ws.Range("I2").Value = "bobina"
$ws.Range("J2").Value = 169
$ws.Range("K2").Value = 860
$ws.Range("L2").Value = "CAMPO VUOTO"
$ws.Range("M2").Value = "Dati OK"
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = "CAMPO VUOTO"
$ws.Range("P2").Value = "CAMPO VUOTO"
$ws.Range("A3").Value = 251702
$ws.Range("B3").Value = 45903
$ws.Range("C3").Value = "STAMPATO"
$ws.Range("D3").Value = 45853
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 60063
$ws.Range("G3").Value = 4607
$ws.Range("H3").Value = "1"
$ws.Range("I3").Value = "bobina"
$ws.Range("J3").Value = 340
$ws.Range("K3").Value = 1020
$ws.Range("L3").Value = "CAMPO VUOTO"
$ws.Range("M3").Value = "Dati OK"
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = "CAMPO VUOTO"
$ws.Range("P3").Value = "CAMPO VUOTO"
$ws.Range("A4").Value = 251310
$ws.Range("B4").Value = 45770
$ws.Range("C4").Value = "STAMPATO"
$ws.Range("D4").Value = 45769
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 8611
$ws.Range("G4").Value = 336
$ws.Range("H4").Value = "6"
$ws.Range("I4").Value = "bobina"
$ws.Range("J4").Value = 410
$ws.Range("K4").Value = 820
$ws.Range("L4").Value = "CAMPO VUOTO"
$ws.Range("M4").Value = "Dati OK"
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = "CAMPO VUOTO"
$ws.Range("P4").Value = "CAMPO VUOTO"
$ws.Range("A5").Value = 252683
$ws.Range("B5").Value = 45903
$ws.Range("C5").Value = "STAMPATO"
$ws.Range("D5").Value = 45901
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 19381
$ws.Range("G5").Value = 2000
$ws.Range("H5").Value = "5"
$ws.Range("I5").Value = "bobina"
$ws.Range("J5").Value = 820
$ws.Range("K5").Value = 820
$ws.Range("L5").Value = "CAMPO VUOTO"
$ws.Range("M5").Value = "Dati OK"
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = "X"
$ws.Range("P5").Value = 40295
$ws.Range("A6").Value = 252459
$ws.Range("B6").Value = 45855
$ws.Range("C6").Value = "STAMPATO"
$ws.Range("D6").Value = "CAMPO VUOTO"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 4000
$ws.Range("G6").Value = 91
$ws.Range("H6").Value = "CAMPO VUOTO"
$ws.Range("I6").Value = "foglio"
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 480
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = "Dati OK"
$ws.Range("N6").Value = 3
$ws.Range("O6").Value = "CAMPO VUOTO"
$ws.Range("P6").Value = "CAMPO VUOTO"
$ws.Range("A7").Value = 252596
$ws.Range("B7").Value = 45904
$ws.Range("C7").Value = "STAMPATO"
$ws.Range("D7").Value = 45873
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 3497
$ws.Range("G7").Value = 305
$ws.Range("H7").Value = "5 / 9"
$ws.Range("I7").Value = "bobina"
$ws.Range("J7").Value = 470
$ws.Range("K7").Value = 940
$ws.Range("L7").Value = "CAMPO VUOTO"
$ws.Range("M7").Value = "Dati OK"
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = "X"
$ws.Range("P7").Value = 40295
$ws.Range("A8").Value = 251624
$ws.Range("B8").Value = 45908
$ws.Range("C8").Value = "IN STAMPA"
$ws.Range("D8").Value = 45895
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 10957
$ws.Range("G8").Value = 574
$ws.Range("H8").Value = "CAMPO VUOTO"
$ws.Range("I8").Value = "bobina"
$ws.Range("J8").Value = 155
$ws.Range("K8").Value = 620
$ws.Range("L8").Value = "CAMPO VUOTO"
$ws.Range("M8").Value = "Dati OK"
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = "X"
$ws.Range("P8").Value = 40308
